$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add new columns P (14) and Q (15), copying the header style from O1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Cells.Item(1, 16).Value = 14
$ws.Cells.Item(1, 17).Value = 15

# New values for columns B..Q (2..17) across rows 2..25
    $data = @(
        @(2.822749202295086,0.9614788761574573,0.2767175078565032,0,1.408598306637757,0.9223593894628266,0.0007673538630890153,0.002592180409469069,0.5871337558215117,0.4126316076062935,0.5212797041709507,0,0,0,0.7135351973630151,0),
        @(2.455957239887937,0.8356714250496111,0.245392408638196,0,1.286228922624701,0.8379586001580321,0.00008923178587139091,0.001313198610051636,0.5521700609348983,0.3965676120329391,0.459498197322489,0,0,0,0.7342630517728814,0),
        @(2.230281723343523,0.7589633814358194,0.2262479012355527,0,1.211850600625056,0.7868067739616151,0.000004634781228141094,0.0008073644887036657,0.531123105327282,0.3870005918597101,0.4216379631896672,0,0,0,0.7476647841110626,0),
        @(2.136675398094894,0.7285752984999476,0.2186339207288768,0,1.180717436397273,0.7651626668646401,0.00004896367825546477,0.0007281875672884297,0.5221619973232094,0.3826128747642734,0.4062030907704468,0,0,0,0.7535138135668795,0),
        @(2.119279639592264,0.7244588001216812,0.2175775413699341,0,1.174344704005307,0.7604170376747561,0.00006130149985716926,0.0008014915406961975,0.5200910119324931,0.381207408937442,0.4036204088172042,0,0,0,0.7547714420188427,0),
        @(2.22396407943387,0.7610794155898475,0.226711335643941,0,1.208113014236986,0.7833397433924745,0.000005180087725875637,0.001015392295610873,0.5293926054184652,0.3850796182895664,0.4213746578916044,0,0,0,0.7485143899553677,0),
        @(2.689590676791681,0.9213498455677893,0.2666507166972281,0,1.361824054054438,0.8888833852381737,0.0004585153634280026,0.002333382612873969,0.5728500880177307,0.4045591808146511,0.4998825134194647,0,0,0,0.7216323048162252,0),
        @(3.609147961137751,1.23618451827997,0.3448369386455283,0,1.67632278748728,1.107835572763733,0.004219505657429456,0.007067408212605208,0.6647295746398925,0.4487607231743524,0.6549394197761842,0,0,0,0.6730455217945135,0),
        @(4.263328620389132,1.461990571264039,0.3916882368237964,0,1.881853825528253,1.255432182840394,0.008889636139635115,0.01233835052667764,0.7260524322674371,0.4761797966285428,0.7362239310369745,0,0,0,0.6454783881740056,0),
        @(4.388964532383795,1.472732244771294,0.3183908058273346,0,1.717637719934785,1.174059843161785,0.02705849694661211,0.01422412188601996,0.6806921604905085,0.4334744804278543,0.510024901169686,0,0,0,0.6733052691925607,0),
        @(4.362738846168327,1.429460209367221,0.2536084389801658,0,1.546738764331565,1.081056140663847,0.06485519204617418,0.01429583449481964,0.6326169235091328,0.3957044321000609,0.3431980961185133,0,0,0,0.7062295608777234,0),
        @(4.21036241518749,1.346575809612489,0.1928830503427719,0,1.356541292745945,0.9690630223347227,0.1193081247849079,0.01324370443371503,0.5772538710299386,0.3566512627693541,0.2130058178676393,0,0,0,0.745296851824115,0),
        @(4.045911200991725,1.272030061168323,0.1537579186106797,0,1.216120358038268,0.8824100458245425,0.1676444444873368,0.01212421460037749,0.5352995972059063,0.3287111816219692,0.1446390055678606,0,0,0,0.7760605669410268,0),
        @(3.977553790582419,1.24660775365885,0.1442963149526406,0,1.176366418865371,0.8561780972982547,0.1797667073301454,0.01172189913482491,0.5230081306106626,0.320971169381373,0.1306222439530735,0,0,0,0.7847038889221807,0),
        @(3.727842765043533,1.169958794965112,0.138692231350575,0,1.129501440597352,0.8158730886115535,0.1654513300767491,0.009739544730046568,0.5078293786882995,0.3165429365586157,0.1270696677652836,0,0,0,0.7860759270882767,0),
        @(3.623193811179931,1.150597803977064,0.1541821374168819,0,1.166659342950268,0.8303202184770413,0.1269825198741756,0.008845415174827842,0.51821517056554,0.3277339513760253,0.1561604325838886,0,0,0,0.7710341479843734,0),
        @(3.634819605601194,1.176251198535624,0.1933310325455295,0,1.286224895920171,0.8970455888971145,0.07422037896036926,0.00845816410150313,0.5534210837172111,0.3555926074306939,0.2370361689147273,0,0,0,0.7410236455050878,0),
        @(3.72961766710489,1.241752552670278,0.2558834251305484,0,1.46510800576354,0.9986822874588199,0.02988699570717301,0.008989447562928632,0.60463624880569,0.3932869801745582,0.3849186904978978,0,0,0,0.7069972591236535,0),
        @(4.074419593190726,1.410015127782628,0.3806429557066764,0,1.815779375810862,1.205572243439804,0.007462225637297148,0.01142981867727499,0.7043488746342916,0.4627705454633428,0.7132122737094164,0,0,0,0.6555849926175057,0),
        @(4.608916836899198,1.600606120937471,0.4346280157941749,0,2.024122496721318,1.349100121162508,0.01209883623861852,0.01623267062242828,0.7664489851055407,0.4954464898680868,0.8284808207118886,0,0,0,0.6293617388265389,0),
        @(4.956165040578583,1.717335136006739,0.4634353574261922,0,2.15119711269729,1.439978334399569,0.01554325234775444,0.01951508287095027,0.8056823073793282,0.5166188975351531,0.887124089543633,0,0,0,0.613715066504831,0),
        @(4.777025314345906,1.651831134489271,0.4473474894068374,0,2.087199253803846,1.395151579111911,0.01367040244171502,0.01748884168139853,0.78658652420404,0.5074707950997777,0.8558543570032242,0,0,0,0.6208194040438713,0),
        @(4.089903469183014,1.411021295005412,0.3879680861972048,0,1.841999468580966,1.222386791923185,0.007558988817717049,0.01107136714276535,0.7128368832806018,0.4705109541139336,0.7381943719009314,0,0,0,0.6512588292133472,0),
        @(3.351137541555147,1.155000107763783,0.3246061826837376,0,1.584260544369428,1.041990055597253,0.002850567825891659,0.005831604424576753,0.6364806635106106,0.4330984824503759,0.6127187049490885,0,0,0,0.6871700041167728,0)
    )


for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $rowNum = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $colNum = $j + 2
        $ws.Cells.Item($rowNum, $colNum).Value = $row[$j]
    }
}
